$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A45').Value = '44. (Passage Reading: Azerbaijan’s transformation is particularly impressive considering that twenty-seven years ago it was involved in a deep political and economic crisis. Amid post-Soviet uncertainty, the economy was in free fall, with a severe contraction in output and growing inflation. The revival of the country’s economy followed the introduction of a new oil strategy which encompassed a series of production sharing agreements, and the opening of Azerbaijan’s historic oil sector to international investment. A 1994 agreement, hailed as the ‘contract of the century’, gave investors access to Azerbaijan’s rich oil fields in the Caspian Sea. The resulting inflow of capital and expertise permanently altered the country’s fortunes. Large-scale production and the construction of multiple pipelines contributed to making Azerbaijan one of the world’s largest oil and gas exporters. Azerbaijan’s biggest challenge today is diversifying away from hydrocarbons (chemicals found in petrol and natural gas), which represent more than 90 percent of the country’s total exports, a 60 percent increase from fifteen years ago. The country expects to capitalise on hydrocarbons for at least another forty years, but oil production is believed to have peaked and Azerbaijan must do more for its non-oil economy. For the government this is a clear priority.) It can be understood from the passage that ----.'
$ws.Range('A46').Value = '45. (Passage Reading: Azerbaijan’s transformation is particularly impressive considering that twenty-seven years ago it was involved in a deep political and economic crisis. Amid post-Soviet uncertainty, the economy was in free fall, with a severe contraction in output and growing inflation. The revival of the country’s economy followed the introduction of a new oil strategy which encompassed a series of production sharing agreements, and the opening of Azerbaijan’s historic oil sector to international investment. A 1994 agreement, hailed as the ‘contract of the century’, gave investors access to Azerbaijan’s rich oil fields in the Caspian Sea. The resulting inflow of capital and expertise permanently altered the country’s fortunes. Large-scale production and the construction of multiple pipelines contributed to making Azerbaijan one of the world’s largest oil and gas exporters. Azerbaijan’s biggest challenge today is diversifying away from hydrocarbons (chemicals found in petrol and natural gas), which represent more than 90 percent of the country’s total exports, a 60 percent increase from fifteen years ago. The country expects to capitalise on hydrocarbons for at least another forty years, but oil production is believed to have peaked and Azerbaijan must do more for its non-oil economy. For the government this is a clear priority.) The underlined word in the passage ‘hailed’ is closest in meaning to ----.'
$ws.Range('A47').Value = '46. (Passage Reading: Azerbaijan’s transformation is particularly impressive considering that twenty-seven years ago it was involved in a deep political and economic crisis. Amid post-Soviet uncertainty, the economy was in free fall, with a severe contraction in output and growing inflation. The revival of the country’s economy followed the introduction of a new oil strategy which encompassed a series of production sharing agreements, and the opening of Azerbaijan’s historic oil sector to international investment. A 1994 agreement, hailed as the ‘contract of the century’, gave investors access to Azerbaijan’s rich oil fields in the Caspian Sea. The resulting inflow of capital and expertise permanently altered the country’s fortunes. Large-scale production and the construction of multiple pipelines contributed to making Azerbaijan one of the world’s largest oil and gas exporters. Azerbaijan’s biggest challenge today is diversifying away from hydrocarbons (chemicals found in petrol and natural gas), which represent more than 90 percent of the country’s total exports, a 60 percent increase from fifteen years ago. The country expects to capitalise on hydrocarbons for at least another forty years, but oil production is believed to have peaked and Azerbaijan must do more for its non-oil economy. For the government this is a clear priority.) According to the passage, for further economic development Azerbaijan needs to ----.'
$ws.Range('A49').Value = '48. (Passage Reading: Deconstruction – the method which claims there is no single explanation of the meaning of a piece of writing – in and of itself, is almost dead, but it retains enormous influence on current critical thinking. While deconstruction itself peaked in the late 1970s, an evolving version of deconstruction has been absorbed by later thinking and often remains crucial to, and even taken for granted by, contemporary cultural and literary criticism. The founding figure and the intellectual force behind deconstruction was Jacques Derrida. His first widely influential book, Of Grammatology, appeared in French in 1967, and many of his later writings recast or extend the ideas of Of Grammatology for other contexts and issues. Advocates and popularisers saw his ideas as so revolutionary that they put him in the same light as Copernicus, Newton, Einstein, and Freud. Although Derrida was primarily a philosopher, his work had little impact on philosophy until after its shock wave hit literary studies. Deconstruction has since carried its influence far beyond literary studies, not only to philosophy but also to all of the humanities and often to the social sciences.) It is stated in the passage that Jacques Derrida’s book Of Grammatology ----.'
$ws.Range('A50').Value = '49.(Passage Reading: Deconstruction – the method which claims there is no single explanation of the meaning of a piece of writing – in and of itself, is almost dead, but it retains enormous influence on current critical thinking. While deconstruction itself peaked in the late 1970s, an evolving version of deconstruction has been absorbed by later thinking and often remains crucial to, and even taken for granted by, contemporary cultural and literary criticism. The founding figure and the intellectual force behind deconstruction was Jacques Derrida. His first widely influential book, Of Grammatology, appeared in French in 1967, and many of his later writings recast or extend the ideas of Of Grammatology for other contexts and issues. Advocates and popularisers saw his ideas as so revolutionary that they put him in the same light as Copernicus, Newton, Einstein, and Freud. Although Derrida was primarily a philosopher, his work had little impact on philosophy until after its shock wave hit literary studies. Deconstruction has since carried its influence far beyond literary studies, not only to philosophy but also to all of the humanities and often to the social sciences.) Derrida is seen in the same line with Copernicus, Newton, Einstein, and Freud because he ----.'
$ws.Range('A51').Value = '50. (Passage Reading: Deconstruction – the method which claims there is no single explanation of the meaning of a piece of writing – in and of itself, is almost dead, but it retains enormous influence on current critical thinking. While deconstruction itself peaked in the late 1970s, an evolving version of deconstruction has been absorbed by later thinking and often remains crucial to, and even taken for granted by, contemporary cultural and literary criticism. The founding figure and the intellectual force behind deconstruction was Jacques Derrida. His first widely influential book, Of Grammatology, appeared in French in 1967, and many of his later writings recast or extend the ideas of Of Grammatology for other contexts and issues. Advocates and popularisers saw his ideas as so revolutionary that they put him in the same light as Copernicus, Newton, Einstein, and Freud. Although Derrida was primarily a philosopher, his work had little impact on philosophy until after its shock wave hit literary studies. Deconstruction has since carried its influence far beyond literary studies, not only to philosophy but also to all of the humanities and often to the social sciences.) According to the passage, which of the following is true about deconstruction?'
$ws.Range('A53').Value = '52. (Passage Reading: What goes around comes around. Our greenhouse gas emissions do not just warm the planet, they also acidify the oceans; however, it is a fact that most people are not aware of. Now it turns out that the changes in ocean chemistry they cause will feed back into the climate, further driving up temperatures. This was not always thought to be true. Climatologists had considered that carbon dioxide (CO2) absorbed by the ocean was stored, and unable to affect the climate. But a study now suggests that the acidification it causes will rebound on the entire planet, by acting on tiny marine plants called phytoplankton. These plants produce a chemical called dimethyl sulphide (DMS) that drifts up into the air and reflects sunlight back into space, cooling the planet. Katharina Six and her team at the Max Planck Institute for Meteorology in Hamburg gathered experimental data showing that phytoplankton produce less DMS as seawater becomes less alkaline due to acidification. After entering these figures into climate models, they estimate that 18 percent less DMS will be released from the oceans by the year 2100, compared to pre-industrial times. If the concentration of CO2 in the atmosphere doubles, which is likely to happen later this century, temperatures are expected to rise between 2 and 4.5°C. Ocean acidification will add between 0.23 and 0.48°C to that figure, Six estimates.)   It is understood from the passage that the chemical dimethyl sulphide ----.'
$ws.Range('A54').Value = '53. (Passage Reading: What goes around comes around. Our greenhouse gas emissions do not just warm the planet, they also acidify the oceans; however, it is a fact that most people are not aware of. Now it turns out that the changes in ocean chemistry they cause will feed back into the climate, further driving up temperatures. This was not always thought to be true. Climatologists had considered that carbon dioxide (CO2) absorbed by the ocean was stored, and unable to affect the climate. But a study now suggests that the acidification it causes will rebound on the entire planet, by acting on tiny marine plants called phytoplankton. These plants produce a chemical called dimethyl sulphide (DMS) that drifts up into the air and reflects sunlight back into space, cooling the planet. Katharina Six and her team at the Max Planck Institute for Meteorology in Hamburg gathered experimental data showing that phytoplankton produce less DMS as seawater becomes less alkaline due to acidification. After entering these figures into climate models, they estimate that 18 percent less DMS will be released from the oceans by the year 2100, compared to pre-industrial times. If the concentration of CO2 in the atmosphere doubles, which is likely to happen later this century, temperatures are expected to rise between 2 and 4.5°C. Ocean acidification will add between 0.23 and 0.48°C to that figure, Six estimates.)   According to the data collected by Six and her team, ----.'
$ws.Range('A55').Value = '54. (Passage Reading: What goes around comes around. Our greenhouse gas emissions do not just warm the planet, they also acidify the oceans; however, it is a fact that most people are not aware of. Now it turns out that the changes in ocean chemistry they cause will feed back into the climate, further driving up temperatures. This was not always thought to be true. Climatologists had considered that carbon dioxide (CO2) absorbed by the ocean was stored, and unable to affect the climate. But a study now suggests that the acidification it causes will rebound on the entire planet, by acting on tiny marine plants called phytoplankton. These plants produce a chemical called dimethyl sulphide (DMS) that drifts up into the air and reflects sunlight back into space, cooling the planet. Katharina Six and her team at the Max Planck Institute for Meteorology in Hamburg gathered experimental data showing that phytoplankton produce less DMS as seawater becomes less alkaline due to acidification. After entering these figures into climate models, they estimate that 18 percent less DMS will be released from the oceans by the year 2100, compared to pre-industrial times. If the concentration of CO2 in the atmosphere doubles, which is likely to happen later this century, temperatures are expected to rise between 2 and 4.5°C. Ocean acidification will add between 0.23 and 0.48°C to that figure, Six estimates.)   What is the main aim of the author of this passage?'
$ws.Range('A57').Value = '56. (Passage Reading: So far, only two types of mental exertion have been shown to improve or preserve overall cognitive ability. One is musical training. Experts state that listening to and making music is not only an auditory experience; it is a multisensory and motor experience. Playing an instrument changes how the brain interprets and integrates a wide range of sensory information, and making music over a long period of time can change brain function and brain structure. The other is to learn a new language. As for the role of language-learning in brain power, a study carried out among students revealed that attention gets better when a new language is studied intensively even if it is for a short time and it gets even better when more practice is done for longer. The results of the study confirm the cognitive benefits of language-learning and give important and warming messages: firstly, even a short, intensive course can show beneficial effects on some cognitive functions. Secondly, this effect can be maintained through practice. So, we can say that it is never too late to start a mental activity such as learning a new language.) It is pointed out in the passage that learning a new language ----.'
$ws.Range('A58').Value = '57. (Passage Reading: So far, only two types of mental exertion have been shown to improve or preserve overall cognitive ability. One is musical training. Experts state that listening to and making music is not only an auditory experience; it is a multisensory and motor experience. Playing an instrument changes how the brain interprets and integrates a wide range of sensory information, and making music over a long period of time can change brain function and brain structure. The other is to learn a new language. As for the role of language-learning in brain power, a study carried out among students revealed that attention gets better when a new language is studied intensively even if it is for a short time and it gets even better when more practice is done for longer. The results of the study confirm the cognitive benefits of language-learning and give important and warming messages: firstly, even a short, intensive course can show beneficial effects on some cognitive functions. Secondly, this effect can be maintained through practice. So, we can say that it is never too late to start a mental activity such as learning a new language.) According to the passage, the study done on students shows that ----.'
$ws.Range('A59').Value = '58.(Passage Reading: So far, only two types of mental exertion have been shown to improve or preserve overall cognitive ability. One is musical training. Experts state that listening to and making music is not only an auditory experience; it is a multisensory and motor experience. Playing an instrument changes how the brain interprets and integrates a wide range of sensory information, and making music over a long period of time can change brain function and brain structure. The other is to learn a new language. As for the role of language-learning in brain power, a study carried out among students revealed that attention gets better when a new language is studied intensively even if it is for a short time and it gets even better when more practice is done for longer. The results of the study confirm the cognitive benefits of language-learning and give important and warming messages: firstly, even a short, intensive course can show beneficial effects on some cognitive functions. Secondly, this effect can be maintained through practice. So, we can say that it is never too late to start a mental activity such as learning a new language.) Which of the following best describes the author’s attitude towards learning a new language?'
$ws.Range('A61').Value = '60. (Passage Reading: So far, only two types of mental exertion have been shown to improve or preserve overall cognitive ability. One is musical training. Experts state that listening to and making music is not only an auditory experience; it is a multisensory and motor experience. Playing an instrument changes how the brain interprets and integrates a wide range of sensory information, and making music over a long period of time can change brain function and brain structure. The other is to learn a new language. As for the role of language-learning in brain power, a study carried out among students revealed that attention gets better when a new language is studied intensively even if it is for a short time and it gets even better when more practice is done for longer. The results of the study confirm the cognitive benefits of language-learning and give important and warming messages: firstly, even a short, intensive course can show beneficial effects on some cognitive functions. Secondly, this effect can be maintained through practice. So, we can say that it is never too late to start a mental activity such as learning a new language.) What point does Richmond make about contamination in the passage?'
$ws.Range('A63').Value = '62. (Passage Reading: When we take medication, some can pass through us intact and go down the toilet. Once in the sewers, these drugs are not removed during the treatment process, and end up in the waterways. To find out how pharmaceutical waste finds its way into aquatic creatures, Erinn Richmond at Monash University in Australia and her colleagues sampled flies, beetles, other insects and spiders from six waterways in the greater Melbourne region. Overall, the group detected 69 medications in the invertebrates. Organisms near wastewater plants contained the highest levels, but low levels were also detected in those from more pristine areas. Richmond presumes river-borne pharmaceuticals probably accumulate in flies and beetles while they are underwater larvae, then transfer to spiders that feed on them after they emerge as adults. Other predators in the food web such as fish, platypuses, birds, bats and frogs may also become cross-contaminated, she says. Her team calculated that trout and platypuses feeding on insects near wastewater plants would consume anti-depressant doses up to half the level normally prescribed to humans. The impacts of this are not certain, but previous studies have found that high-level exposure to anti-depressants causes fish to become less vigilant about avoiding predators. Male birds also seem to sing less to females when exposed to the drugs.) What could be the best title for this passage?'
